# BWP Object Repository and fixes for VRelay
# Updates the "DateProd" (last-run timestamp) column on the generic VT test
# sheets to reflect a new Katalon run, and flips one ResultProd cell from
# Pass to Fail (VT-P-DebitCredit-DualCF-Generic, row 3).

$wb = $excel.ActiveWorkbook

function Set-Cell($sheetName, $cellRef, $value) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range($cellRef).Value = $value
}

# VT-P-DebitVoid-DualCF-Generic
Set-Cell "VT-P-DebitVoid-DualCF-Generic" "B2" "Thu Sep 04 07:17:22 IST 2025"

# VT-P-DebitVoid-SingleCF-Generic
Set-Cell "VT-P-DebitVoid-SingleCF-Generic" "B2" "Thu Sep 04 07:23:19 IST 2025"
Set-Cell "VT-P-DebitVoid-SingleCF-Generic" "B3" "Thu Sep 04 07:24:25 IST 2025"
Set-Cell "VT-P-DebitVoid-SingleCF-Generic" "B4" "Thu Sep 04 07:25:35 IST 2025"
Set-Cell "VT-P-DebitVoid-SingleCF-Generic" "B5" "Thu Sep 04 07:26:41 IST 2025"

# VT-P-DebitVoid-NoCF-Generic
Set-Cell "VT-P-DebitVoid-NoCF-Generic" "B2" "Thu Sep 04 07:18:37 IST 2025"
Set-Cell "VT-P-DebitVoid-NoCF-Generic" "B3" "Thu Sep 04 07:19:43 IST 2025"
Set-Cell "VT-P-DebitVoid-NoCF-Generic" "B4" "Thu Sep 04 07:21:01 IST 2025"
Set-Cell "VT-P-DebitVoid-NoCF-Generic" "B5" "Thu Sep 04 07:22:09 IST 2025"

# VT-P-DebitCredit-DualCF-Generic (row 3 also flips ResultProd Pass -> Fail)
Set-Cell "VT-P-DebitCredit-DualCF-Generic" "B2" "Thu Sep 04 07:03:43 IST 2025"
Set-Cell "VT-P-DebitCredit-DualCF-Generic" "A3" "Fail"
Set-Cell "VT-P-DebitCredit-DualCF-Generic" "B3" "Thu Sep 04 07:05:07 IST 2025"
Set-Cell "VT-P-DebitCredit-DualCF-Generic" "B4" "Thu Sep 04 07:06:15 IST 2025"
Set-Cell "VT-P-DebitCredit-DualCF-Generic" "B5" "Thu Sep 04 07:07:22 IST 2025"

# VT-P-DebitCredit-SingleCF-Gener
Set-Cell "VT-P-DebitCredit-SingleCF-Gener" "B2" "Thu Sep 04 07:12:50 IST 2025"
Set-Cell "VT-P-DebitCredit-SingleCF-Gener" "B3" "Thu Sep 04 07:14:03 IST 2025"
Set-Cell "VT-P-DebitCredit-SingleCF-Gener" "B4" "Thu Sep 04 07:15:12 IST 2025"
Set-Cell "VT-P-DebitCredit-SingleCF-Gener" "B5" "Thu Sep 04 07:16:20 IST 2025"

# VT-P-DebitCredit-NoCF-Generic
Set-Cell "VT-P-DebitCredit-NoCF-Generic" "B2" "Thu Sep 04 07:08:26 IST 2025"
Set-Cell "VT-P-DebitCredit-NoCF-Generic" "B3" "Thu Sep 04 07:09:33 IST 2025"
Set-Cell "VT-P-DebitCredit-NoCF-Generic" "B4" "Thu Sep 04 07:10:42 IST 2025"
Set-Cell "VT-P-DebitCredit-NoCF-Generic" "B5" "Thu Sep 04 07:11:46 IST 2025"

# VT-C-DebitCredit-DualCF-Generic
Set-Cell "VT-C-DebitCredit-DualCF-Generic" "B2" "Thu Sep 04 07:00:23 IST 2025"

# VT-C-DebitCredit-SingleCF-Gener
Set-Cell "VT-C-DebitCredit-SingleCF-Gener" "B2" "Thu Sep 04 07:02:30 IST 2025"

# VT-C-DebitCredit-NoCF-Generic
Set-Cell "VT-C-DebitCredit-NoCF-Generic" "B2" "Thu Sep 04 07:01:30 IST 2025"
